# Update odds values in row 2, row 3 and row 4 of Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 updates
$ws.Range("G2").Value = 2.5
$ws.Range("I2").Value = 3.2
$ws.Range("J2").Value = 3.5
$ws.Range("M2").Value = 1.13
$ws.Range("N2").Value = 6
$ws.Range("Q2").Value = 2.88
$ws.Range("R2").Value = 1.4
$ws.Range("W2").Value = 5.5
$ws.Range("X2").Value = 10
$ws.Range("Z2").Value = 26
$ws.Range("AA2").Value = 26
$ws.Range("AJ2").Value = 13
$ws.Range("AN2").Value = 4.33
$ws.Range("AW2").Value = 5
$ws.Range("AX2").Value = 21
$ws.Range("AY2").Value = 41
$ws.Range("AZ2").Value = 81
$ws.Range("BB2").Value = 451

# Row 3 updates
$ws.Range("Q3").Value = 2.15
$ws.Range("R3").Value = 1.67

# Row 4 updates
$ws.Range("G4").Value = 2
